# ST1_Line_Setup.xlsx - update row 2 data entry
# - Date/shift changed, several numeric-looking readings reset to 0.1,
#   and the three sign-off cells shortened to "o".
# Values in this template are stored as *text* (not numbers), so for the
# cells whose new value looks like a number we force a text number format
# before writing, then drop back to the Normal style so no stray
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Plain text updates (not numeric-looking, no special handling needed)
$ws.Range("A2").Value = "2025-04-03T15:03"
$ws.Range("B2").Value = "SHIFT_2"

# Numeric-looking text values -> keep stored as text
Set-TextValue "C2" "0.1"
Set-TextValue "D2" "0.1"
Set-TextValue "E2" "0.1"
Set-TextValue "F2" "0.1"

Set-TextValue "M2" "0.1"
Set-TextValue "N2" "0.1"
Set-TextValue "O2" "0.1"
Set-TextValue "P2" "0.1"
Set-TextValue "Q2" "0.1"
Set-TextValue "R2" "0.1"
Set-TextValue "S2" "0.1"
Set-TextValue "T2" "0.1"

# Sign-off cells shortened to "o"
$ws.Range("W2").Value = "o"
$ws.Range("X2").Value = "o"
$ws.Range("Y2").Value = "o"
